$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a brand-new paragraph (about thyristors) right before the
#    "For another version of 12-pulse rectifier..." paragraph.
# ------------------------------------------------------------------
$target = $d.Paragraphs(4)
$target.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs(4)
$newPara.Alignment = 3   # wdAlignParagraphJustify

$nr = $newPara.Range
$nr.Collapse(1)
$nr.InsertAfter("     ")
$nr.Collapse(0)
$nr.InsertAfter(" ")
$nr.Collapse(0)
$nr.InsertAfter("For a different variation in the 12-pulse rectifiers, thyristors can be used instead of diodes.")
$nr.Collapse(0)
$nr.InsertAfter(" ")
$nr.Collapse(0)
$nr.InsertAfter("When thyristors are used, average voltage of output will decrease ")
$nr.Collapse(0)
$nr.InsertAfter("depending on firing angle")
$nr.Collapse(0)
$nr.InsertAfter(" ")
$nr.Collapse(0)
$nr.InsertAfter("compared to ")
$nr.Collapse(0)
$nr.InsertAfter("using")
$nr.Collapse(0)
$nr.InsertAfter(" diodes.Morover,")
$nr.Collapse(0)
$nr.InsertAfter(" ")
$nr.Collapse(0)
$nr.InsertAfter("t")
$nr.Collapse(0)
$nr.InsertAfter("hyristor rectifiers are partially controlled")
$nr.Collapse(0)
$nr.InsertAfter(",whereas d")
$nr.Collapse(0)
$nr.InsertAfter("iodes are not controlled")
$nr.Collapse(0)
$nr.InsertAfter(".")
$nr.Collapse(0)
$nr.InsertAfter(" ")
$nr.Collapse(0)
$nr.InsertAfter("We can control output voltage")
$nr.Collapse(0)
$nr.InsertAfter(" by just changing the firing angle")
$nr.Collapse(0)
$nr.InsertAfter(" of thristor. ")
$nr.Collapse(0)
$nr.InsertAfter("Thus,")
$nr.Collapse(0)
$nr.InsertAfter(" ")
$nr.Collapse(0)
$nr.InsertAfter("thyristors are")
$nr.Collapse(0)
$nr.InsertAfter(" good device")
$nr.Collapse(0)
$nr.InsertAfter("s")
$nr.Collapse(0)
$nr.InsertAfter(" for controlling purpose")
$nr.Collapse(0)
$nr.InsertAfter(".")
$nr.Collapse(0)
$nr.InsertAfter(" ")

# ------------------------------------------------------------------
# 2) Update the "For another version..." paragraph itself (now #5):
#    - justify alignment
#    - extra leading whitespace
#    - ",wye" -> ", wye"  (and split the dash/c runs)
#    - append the new sentence about the phase shift, inserted so the
#      hidden _GoBack bookmark ends up inside the new "connection" word
# ------------------------------------------------------------------
$p2 = $d.Paragraphs(5)
$p2.Alignment = 3   # wdAlignParagraphJustify

$p2.Range.Find.Execute(" F", $true, $false, $false, $false, $false, $true, 1, $false, "     F", 2)

$p2.Range.Find.Execute("onnection is used in primary side.", $true, $false, $false, $false, $false, $true, 1, $false, "tion) and secondary side(delta connection) by 30°.", 2)

$p2.Range.Find.Execute(",wye-", $true, $false, $false, $false, $false, $true, 1, $false, ", wye-connection is used in primary side.There is a phase shift between primary side(wye conne", 2)

Write-Output "done"
